# Update "想去人数" (want-to-go count) figures that changed between crawls.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 620
$ws1.Range("F5").Value = 4526
$ws1.Range("F6").Value = 1819
$ws1.Range("F8").Value = 126
$ws1.Range("F9").Value = 3057
$ws1.Range("F12").Value = 239
$ws1.Range("F13").Value = 583
$ws1.Range("F14").Value = 503
$ws1.Range("F15").Value = 505
$ws1.Range("F16").Value = 351
$ws1.Range("F18").Value = 1751
$ws1.Range("F19").Value = 1287
$ws1.Range("F20").Value = 112
$ws1.Range("F21").Value = 1530
$ws1.Range("F22").Value = 123
$ws1.Range("F23").Value = 601
$ws1.Range("F28").Value = 84
$ws1.Range("F29").Value = 117
$ws1.Range("F30").Value = 79
$ws1.Range("F31").Value = 3418
$ws1.Range("F32").Value = 735
$ws1.Range("F34").Value = 224
$ws1.Range("F35").Value = 52
$ws1.Range("F36").Value = 1662

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 33

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 620
$ws4.Range("F5").Value = 4526
$ws4.Range("F6").Value = 1819
$ws4.Range("F8").Value = 126
$ws4.Range("F9").Value = 3057
$ws4.Range("F12").Value = 239
$ws4.Range("F13").Value = 583
$ws4.Range("F14").Value = 503
$ws4.Range("F15").Value = 505
$ws4.Range("F17").Value = 351
$ws4.Range("F19").Value = 1751
$ws4.Range("F20").Value = 1287
$ws4.Range("F21").Value = 112
$ws4.Range("F22").Value = 1530
$ws4.Range("F23").Value = 123
$ws4.Range("F24").Value = 601
$ws4.Range("F29").Value = 84
$ws4.Range("F30").Value = 117
$ws4.Range("F31").Value = 79
$ws4.Range("F32").Value = 3418
$ws4.Range("F33").Value = 33
$ws4.Range("F34").Value = 735
$ws4.Range("F36").Value = 224
$ws4.Range("F37").Value = 52
$ws4.Range("F38").Value = 1662
